{"js": "// Replace the date line and the 25 \"NN\u00d7NN=\" multiplication problems\n// with their updated values, per the authored diff.\nconst replacements = [\n  [\"2025-05-27 Tuesday\", \"2025-05-28 Wednesday\"],\n  [\"68\u00d745=\", \"19\u00d732=\"],\n  [\"60\u00d772=\", \"50\u00d744=\"],\n  [\"77\u00d783=\", \"62\u00d754=\"],\n  [\"97\u00d796=\", \"36\u00d725=\"],\n  [\"86\u00d779=\", \"20\u00d760=\"],\n  [\"41\u00d787=\", \"64\u00d748=\"],\n  [\"76\u00d738=\", \"94\u00d735=\"],\n  [\"27\u00d768=\", \"68\u00d783=\"],\n  [\"11\u00d723=\", \"70\u00d731=\"],\n  [\"36\u00d737=\", \"90\u00d729=\"],\n  [\"56\u00d768=\", \"20\u00d714=\"],\n  [\"68\u00d743=\", \"23\u00d751=\"],\n  [\"41\u00d772=\", \"71\u00d736=\"],\n  [\"30\u00d781=\", \"26\u00d794=\"],\n  [\"20\u00d743=\", \"22\u00d794=\"],\n  [\"39\u00d712=\", \"69\u00d767=\"],\n  [\"64\u00d734=\", \"84\u00d731=\"],\n  [\"61\u00d733=\", \"53\u00d766=\"],\n  [\"58\u00d779=\", \"60\u00d738=\"],\n  [\"47\u00d776=\", \"16\u00d785=\"],\n  [\"20\u00d725=\", \"85\u00d781=\"],\n  [\"70\u00d776=\", \"99\u00d786=\"],\n  [\"68\u00d716=\", \"36\u00d763=\"],\n  [\"16\u00d729=\", \"18\u00d788=\"],\n  [\"16\u00d741=\", \"11\u00d754=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 \"NN\u00d7NN=\" multiplication problems\n# with their updated values, per the authored diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-27 Tuesday\", \"2025-05-28 Wednesday\"),\n    @(\"68\u00d745=\", \"19\u00d732=\"),\n    @(\"60\u00d772=\", \"50\u00d744=\"),\n    @(\"77\u00d783=\", \"62\u00d754=\"),\n    @(\"97\u00d796=\", \"36\u00d725=\"),\n    @(\"86\u00d779=\", \"20\u00d760=\"),\n    @(\"41\u00d787=\", \"64\u00d748=\"),\n    @(\"76\u00d738=\", \"94\u00d735=\"),\n    @(\"27\u00d768=\", \"68\u00d783=\"),\n    @(\"11\u00d723=\", \"70\u00d731=\"),\n    @(\"36\u00d737=\", \"90\u00d729=\"),\n    @(\"56\u00d768=\", \"20\u00d714=\"),\n    @(\"68\u00d743=\", \"23\u00d751=\"),\n    @(\"41\u00d772=\", \"71\u00d736=\"),\n    @(\"30\u00d781=\", \"26\u00d794=\"),\n    @(\"20\u00d743=\", \"22\u00d794=\"),\n    @(\"39\u00d712=\", \"69\u00d767=\"),\n    @(\"64\u00d734=\", \"84\u00d731=\"),\n    @(\"61\u00d733=\", \"53\u00d766=\"),\n    @(\"58\u00d779=\", \"60\u00d738=\"),\n    @(\"47\u00d776=\", \"16\u00d785=\"),\n    @(\"20\u00d725=\", \"85\u00d781=\"),\n    @(\"70\u00d776=\", \"99\u00d786=\"),\n    @(\"68\u00d716=\", \"36\u00d763=\"),\n    @(\"16\u00d729=\", \"18\u00d788=\"),\n    @(\"16\u00d741=\", \"11\u00d754=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
